# Fix IFRS financial figures for 한진칼 (Hanjin KAL) company_list sheet
# The previously entered rows 2-9 had erroneous (inflated) values; replace
# them with the corrected figures, and drop stray U7/U8/U9 cells that no
# longer apply after the correction.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 6250
$ws.Range("E2").Value = 755
$ws.Range("F2").Value = 755
$ws.Range("G2").Value = 3862
$ws.Range("H2").Value = 2993
$ws.Range("I2").Value = 2175
$ws.Range("J2").Value = 818
$ws.Range("K2").Value = 24168
$ws.Range("L2").Value = 7866
$ws.Range("M2").Value = 16302
$ws.Range("N2").Value = 13995
$ws.Range("O2").Value = 2307
$ws.Range("P2").Value = 1325
$ws.Range("Q2").Value = 1178
$ws.Range("R2").Value = -847
$ws.Range("S2").Value = -321
$ws.Range("T2").Value = 1039
$ws.Range("U2").Value = 138
$ws.Range("V2").Value = 4971
$ws.Range("W2").Value = 12.08
$ws.Range("X2").Value = 47.88
$ws.Range("Y2").Value = 21.9
$ws.Range("Z2").Value = 15.58
$ws.Range("AA2").Value = 48.25
$ws.Range("AB2").Value = 988.51
$ws.Range("AC2").Value = 6505
$ws.Range("AD2").Value = 4.62
$ws.Range("AE2").Value = 26490
$ws.Range("AF2").Value = 1.13
$ws.Range("AG2").Value = 72
$ws.Range("AH2").Value = 0.24
$ws.Range("AI2").Value = 1.77
$ws.Range("AJ2").Value = 54307120

# Row 3
$ws.Range("D3").Value = 7223
$ws.Range("E3").Value = 743
$ws.Range("F3").Value = 743
$ws.Range("G3").Value = -1913
$ws.Range("H3").Value = -2052
$ws.Range("I3").Value = -2164
$ws.Range("J3").Value = 112
$ws.Range("K3").Value = 23150
$ws.Range("L3").Value = 8731
$ws.Range("M3").Value = 14419
$ws.Range("N3").Value = 13223
$ws.Range("O3").Value = 1197
$ws.Range("P3").Value = 1333
$ws.Range("Q3").Value = 1194
$ws.Range("R3").Value = -930
$ws.Range("S3").Value = -40
$ws.Range("T3").Value = 72
$ws.Range("U3").Value = 1122
$ws.Range("V3").Value = 5203
$ws.Range("W3").Value = 10.28
$ws.Range("X3").Value = -28.41
$ws.Range("Y3").Value = -15.9
$ws.Range("Z3").Value = -8.67
$ws.Range("AA3").Value = 60.55
$ws.Range("AB3").Value = 816.45
$ws.Range("AC3").Value = -3935
$ws.Range("AD3").Value = -4.85
$ws.Range("AE3").Value = 23981
$ws.Range("AF3").Value = 0.8
$ws.Range("AG3").Value = 72
$ws.Range("AH3").Value = 0.38
$ws.Range("AI3").Value = -1.85
$ws.Range("AJ3").Value = 54607470

# Row 4
$ws.Range("D4").Value = 9910
$ws.Range("E4").Value = 990
$ws.Range("F4").Value = 990
$ws.Range("G4").Value = -4519
$ws.Range("H4").Value = -3964
$ws.Range("I4").Value = -4068
$ws.Range("J4").Value = 104
$ws.Range("K4").Value = 20620
$ws.Range("L4").Value = 9591
$ws.Range("M4").Value = 11030
$ws.Range("N4").Value = 10017
$ws.Range("O4").Value = 1013
$ws.Range("P4").Value = 1493
$ws.Range("Q4").Value = 1388
$ws.Range("R4").Value = -2570
$ws.Range("S4").Value = 807
$ws.Range("T4").Value = 49
$ws.Range("U4").Value = 1339
$ws.Range("V4").Value = 6172
$ws.Range("W4").Value = 9.99
$ws.Range("X4").Value = -40
$ws.Range("Y4").Value = -35.01
$ws.Range("Z4").Value = -18.11
$ws.Range("AA4").Value = 86.95
$ws.Range("AB4").Value = 415.07
$ws.Range("AC4").Value = -7069
$ws.Range("AD4").Value = -2.17
$ws.Range("AE4").Value = 16779
$ws.Range("AF4").Value = 0.91
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 59170458

# Row 5
$ws.Range("D5").Value = 11497
$ws.Range("E5").Value = 1153
$ws.Range("F5").Value = 1153
$ws.Range("G5").Value = 2676
$ws.Range("H5").Value = 2291
$ws.Range("I5").Value = 2219
$ws.Range("J5").Value = 72
$ws.Range("K5").Value = 27758
$ws.Range("L5").Value = 10661
$ws.Range("M5").Value = 17097
$ws.Range("N5").Value = 15214
$ws.Range("O5").Value = 1883
$ws.Range("P5").Value = 1493
$ws.Range("Q5").Value = 1644
$ws.Range("R5").Value = -3130
$ws.Range("S5").Value = 3140
$ws.Range("T5").Value = 169
$ws.Range("U5").Value = 1475
$ws.Range("V5").Value = 6347
$ws.Range("W5").Value = 10.03
$ws.Range("X5").Value = 19.93
$ws.Range("Y5").Value = 17.59
$ws.Range("Z5").Value = 9.47
$ws.Range("AA5").Value = 62.35
$ws.Range("AB5").Value = 567.6
$ws.Range("AC5").Value = 3717
$ws.Range("AD5").Value = 4.92
$ws.Range("AE5").Value = 25484
$ws.Range("AF5").Value = 0.72
$ws.Range("AG5").Value = 125
$ws.Range("AH5").Value = 0.68
$ws.Range("AI5").Value = 3.37
$ws.Range("AJ5").Value = 59170458

# Row 6
$ws.Range("D6").Value = 13049
$ws.Range("E6").Value = 1088
$ws.Range("F6").Value = 1088
$ws.Range("G6").Value = 8
$ws.Range("H6").Value = -177
$ws.Range("I6").Value = -408
$ws.Range("K6").Value = 28113
$ws.Range("L6").Value = 11275
$ws.Range("M6").Value = 16838
$ws.Range("N6").Value = 14818
$ws.Range("P6").Value = 1493
$ws.Range("Q6").Value = 782
$ws.Range("R6").Value = -2824
$ws.Range("S6").Value = 543
$ws.Range("T6").Value = 71
$ws.Range("U6").Value = 710
$ws.Range("V6").Value = 7314
$ws.Range("W6").Value = 8.34
$ws.Range("X6").Value = -1.35
$ws.Range("Y6").Value = -2.72
$ws.Range("Z6").Value = -0.63
$ws.Range("AA6").Value = 66.96
$ws.Range("AB6").Value = 489.94
$ws.Range("AC6").Value = -683
$ws.Range("AD6").Value = -43.61
$ws.Range("AE6").Value = 24821
$ws.Range("AF6").Value = 1.2
$ws.Range("AG6").Value = 300
$ws.Range("AH6").Value = 1.01
$ws.Range("AI6").Value = -43.93
$ws.Range("AJ6").Value = 59170458

# Row 7
$ws.Range("D7").Value = 12303
$ws.Range("E7").Value = 250
$ws.Range("G7").Value = -1910
$ws.Range("H7").Value = -1810
$ws.Range("I7").Value = -2228
$ws.Range("K7").Value = 26800
$ws.Range("L7").Value = 12134
$ws.Range("M7").Value = 14667
$ws.Range("N7").Value = 12262
$ws.Range("P7").Value = 1491
$ws.Range("Q7").Value = 599
$ws.Range("R7").Value = -1249
$ws.Range("S7").Value = -1585
$ws.Range("T7").Value = 216
$ws.Range("W7").Value = 2.03
$ws.Range("X7").Value = -14.71
$ws.Range("Y7").Value = -16.46
$ws.Range("Z7").Value = -6.59
$ws.Range("AA7").Value = 82.73
$ws.Range("AC7").Value = -3732
$ws.Range("AD7").Value = -10.99
$ws.Range("AE7").Value = 20539
$ws.Range("AF7").Value = 2
$ws.Range("AG7").Value = 300
$ws.Range("AH7").Value = 0.73
$ws.Range("AI7").Value = -7.97
$ws.Range("U7").Value = $null

# Row 8
$ws.Range("D8").Value = 12721
$ws.Range("E8").Value = 532
$ws.Range("G8").Value = 423
$ws.Range("H8").Value = 324
$ws.Range("I8").Value = 423
$ws.Range("K8").Value = 27929
$ws.Range("L8").Value = 13311
$ws.Range("M8").Value = 14618
$ws.Range("N8").Value = 12312
$ws.Range("P8").Value = 1491
$ws.Range("Q8").Value = 638
$ws.Range("R8").Value = -1354
$ws.Range("S8").Value = -420
$ws.Range("T8").Value = 328
$ws.Range("W8").Value = 4.18
$ws.Range("X8").Value = 2.55
$ws.Range("Y8").Value = 3.44
$ws.Range("Z8").Value = 1.18
$ws.Range("AA8").Value = 91.06
$ws.Range("AC8").Value = 708
$ws.Range("AD8").Value = 57.92
$ws.Range("AE8").Value = 20623
$ws.Range("AF8").Value = 1.99
$ws.Range("AG8").Value = 300
$ws.Range("AH8").Value = 0.73
$ws.Range("AI8").Value = 42
$ws.Range("U8").Value = $null

# Row 9
$ws.Range("D9").Value = 13820
$ws.Range("E9").Value = 879
$ws.Range("G9").Value = 881
$ws.Range("H9").Value = 675
$ws.Range("I9").Value = 915
$ws.Range("K9").Value = 28178
$ws.Range("L9").Value = 13258
$ws.Range("M9").Value = 14917
$ws.Range("N9").Value = 12851
$ws.Range("P9").Value = 1491
$ws.Range("Q9").Value = 790
$ws.Range("R9").Value = -623
$ws.Range("S9").Value = -457
$ws.Range("T9").Value = 382
$ws.Range("W9").Value = 6.36
$ws.Range("X9").Value = 4.88
$ws.Range("Y9").Value = 7.27
$ws.Range("Z9").Value = 2.41
$ws.Range("AA9").Value = 88.88
$ws.Range("AC9").Value = 1533
$ws.Range("AD9").Value = 26.74
$ws.Range("AE9").Value = 21526
$ws.Range("AF9").Value = 1.9
$ws.Range("AG9").Value = 300
$ws.Range("AH9").Value = 0.73
$ws.Range("AI9").Value = 19.39
$ws.Range("U9").Value = $null
